# Insert a new data row at row 174 (pushing existing rows 174:276 down to 175:277)
# and populate it with the new weekly price entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(174).Insert()

$ws.Range("A174").Value = 10
$ws.Range("B174").Value = "Vega Modelo de Temuco"
$ws.Range("C174").Value = "La Araucanía"
$ws.Range("D174").Value = 44806
$ws.Range("E174").Value = 9
$ws.Range("F174").Value = 100112043
$ws.Range("G174").Value = "Pepino dulce"
$ws.Range("H174").Value = "Cultivar IV Región"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 65
$ws.Range("K174").Value = 19000
$ws.Range("L174").Value = 19000
$ws.Range("M174").Value = 19000
$ws.Range("N174").Value = "`$/bandeja 18 kilos"
$ws.Range("O174").Value = "Provincia de Limarí"
$ws.Range("P174").Value = 1056
$ws.Range("Q174").Value = 18
$ws.Range("R174").Value = "Hortaliza"
